$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.89%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'43.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'6.40%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.576"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-2.34%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08106"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-3.79%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.666"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-1.76%"
$ws.Range("E6").Style = "Normal"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.895"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-4.58%"
$ws.Range("E7").Style = "Normal"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.277"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-5.06%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'2.755"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-6.68%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9370"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.1172"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-6.53%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1897"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.51%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09677"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.96%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.04152"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'4.81%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.1068"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.37%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001277"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-2.79%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.005941"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.88%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.570"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'3.94%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.75%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'8.550"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-6.60%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1362"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-0.21%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2583"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'2.78%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04312"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.21%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001238"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.67%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004394"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.05%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'3.18%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003985"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.34%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02670"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-5.78%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05481"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.64%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.01143"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'27.45%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007691"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-2.63%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1399"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-2.72%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002107"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.13%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009759"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.38%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00007004"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.26%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000749"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.26%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003534"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'10.07%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002268"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.67%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.26%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.26%"
$ws.Range("E51").Style = "Normal"
